$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update odds on row 5 (Sagan Tosu vs Yokohama F. Marinos)
$ws.Range("Q5").Value = 1.36
$ws.Range("R5").Value = 3.1

# Update odds on row 6 (Shonan Bellmare vs Hokkaido Consadole Sapporo)
$ws.Range("N6").Value = 17
$ws.Range("O6").Value = 1.17
$ws.Range("P6").Value = 5
$ws.Range("Q6").Value = 1.57
$ws.Range("R6").Value = 2.35

# Remove the Mazatlan FC vs U.N.A.M.- Pumas row (row 7); the following
# row (Los Angeles FC vs Vancouver Whitecaps) shifts up to become row 7,
# and the table shrinks by one row overall.
$ws.Range("A7:BD7").EntireRow.Delete()
